$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 33: the "is_active" boolean cell (I33) was missing the
# horizontal-left alignment style that every other data row in that column has.
$ws.Range("I33").HorizontalAlignment = -4131   # xlLeft

# --- Append three new user rows (34-36), matching the existing data pattern.
$newRows = @(
    @{ A=110033; B=9317596771; C="Nikola Tesla"; D="nikola.tesla@xyz.com"; E=818876434 },
    @{ A=110034; B=9317596772; C="Graham Bell";  D="graham.bell@xyz.com";  E=818876435 },
    @{ A=110035; B=9317596773; C="Albert Miles"; D="albert.miles@xyz.com"; E=818876436 }
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = "ACT"
    $ws.Cells.Item($r, 7).Value = "eng"
    $ws.Cells.Item($r, 8).Value = "PWD"
    $ws.Cells.Item($r, 9).Value = $true
    $ws.Range("I$r").HorizontalAlignment = -4131   # xlLeft, match other is_active cells
    $ws.Cells.Item($r, 10).Value = "superadmin"
    $ws.Cells.Item($r, 11).Value = "now()"
    $ws.Cells.Item($r, 12).Value = "now()"
    $r = $r + 1
}

# --- Restore the selection to the top-left of the "unused" area (M1),
# matching the original workbook view before it had drifted to M6.
$ws.Range("M1:XFD1048576").Select() | Out-Null
